# Update LDLC prices history:
# Insert a new "history" column before the existing CC column (nom) so that
# a fresh price snapshot (timestamped 2026-01-31 09:17:47) is recorded,
# pushing the previous "nom" (name) and "url_produit" columns one column
# to the right (CC->CD, CD->CE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at CC; existing CC/CD shift to CD/CE.
$ws.Columns("CC:CC").Insert()

# Header for the newly inserted snapshot column.
$ws.Range("CC1").Value = "2026-01-31 09:17:47"

# For every product row, the new snapshot column just duplicates the most
# recent previous snapshot (now in column CB) as the current price.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $prevVal = $ws.Range("CB$r").Value2
    if ($prevVal -ne $null -and $prevVal -ne "") {
        $ws.Range("CC$r").Value = $prevVal
    }
}
